$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# New column widths for the 3 new columns (G, H, I) added to the table.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 43.59
$ws.Columns.Item(8).ColumnWidth = 25.59
$ws.Columns.Item(9).ColumnWidth = 25.25

# ---------------------------------------------------------------------------
# New data written in the exact order needed so that the new shared-string
# table entries line up the same way they do in the target workbook
# (Notification, Notification Audience, Move Classification, NoActionPlacement,
#  Everyone, SinglePaymentRequired.../MultiplePaymentRequired..., GainCards,
#  MultiplePaymentRequired(Total Cost of Card), SinglePaymentRequired(Total
#  Cost of Card), CancelActionPlacement, ForcedStealPlacement,
#  ForcedStealSetPlacement).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Notification"
$ws.Range("I1").Value = "Notification Audience"
$ws.Range("G1").Value = "Move Classification"

$ws.Range("G2").Value = "NoActionPlacement"
$ws.Range("I2").Value = "Everyone"
$ws.Range("G3").Value = "NoActionPlacement"
$ws.Range("I3").Value = "Everyone"
$ws.Range("G4").Value = "NoActionPlacement"
$ws.Range("I4").Value = "Everyone"
$ws.Range("G5").Value = "NoActionPlacement"
$ws.Range("I5").Value = "Everyone"

$ws.Range("G6").Value = "SinglePaymentRequired(Total Cost of Rent)`nMultiplePaymentRequired(Total Cost of Rent)"

$ws.Range("G7").Value = "NoActionPlacement"
$ws.Range("G8").Value = "NoActionPlacement"
$ws.Range("G9").Value = "NoActionPlacement"
$ws.Range("G10").Value = "NoActionPlacement"
$ws.Range("G11").Value = "NoActionPlacement"

$ws.Range("G12").Value = "GainCards"

$ws.Range("G13").Value = "NoActionPlacement"
$ws.Range("G14").Value = "NoActionPlacement"

$ws.Range("G15").Value = "SinglePaymentRequired(Total Cost of Rent)`nMultiplePaymentRequired(Total Cost of Rent)"

$ws.Range("G16").Value = "NoActionPlacement"
$ws.Range("G17").Value = "NoActionPlacement"

$ws.Range("G18").Value = "MultiplePaymentRequired(Total Cost of Card)"

$ws.Range("G19").Value = "NoActionPlacement"
$ws.Range("G20").Value = "NoActionPlacement"

$ws.Range("G21").Value = "SinglePaymentRequired(Total Cost of Card)"

$ws.Range("G22").Value = "NoActionPlacement"
$ws.Range("G23").Value = "NoActionPlacement"
$ws.Range("G24").Value = "NoActionPlacement"
$ws.Range("G25").Value = "NoActionPlacement"

$ws.Range("G26").Value = "CancelActionPlacement"
$ws.Range("G27").Value = "ForcedStealPlacement"

$ws.Range("G28").Value = "NoActionPlacement"
$ws.Range("G29").Value = "NoActionPlacement"

$ws.Range("G30").Value = "ForcedStealPlacement"

$ws.Range("G31").Value = "NoActionPlacement"
$ws.Range("G32").Value = "NoActionPlacement"

$ws.Range("G33").Value = "ForcedStealSetPlacement"

$ws.Range("G34").Value = "NoActionPlacement"
$ws.Range("G35").Value = "NoActionPlacement"

# ---------------------------------------------------------------------------
# Header row style: white text on a black fill, applied to the whole header
# row A1:I1 (new columns included).
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Font.Color = 16777215
$ws.Range("A1:I1").Interior.Color = 0

# ---------------------------------------------------------------------------
# Wrap text on the two-line notification cells (re-uses the pre-existing
# wrap-text style already used elsewhere in the sheet).
# ---------------------------------------------------------------------------
$ws.Range("G6").WrapText = $true
$ws.Range("G15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 28.8

# ---------------------------------------------------------------------------
# Selection moved to D23 (was D12:D13).
# ---------------------------------------------------------------------------
$ws.Range("D23").Select()
